$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Remove the old "PerfSerach_Classification_TenfoldCV_SMOTE217" row (row 2)
#    This also removes the "2000-2500/10" value that only lived there.
# ---------------------------------------------------------------------------
$ws.Rows("2:2").Delete()

# ---------------------------------------------------------------------------
# 2. Collapse one of the two blank separator rows (old rows 4 & 5, now 3 & 4)
#    so that after we insert the new "Regression" row below, only a single
#    blank row remains between it and the "Ahmed et al." row.
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Delete()

# ---------------------------------------------------------------------------
# 3. Insert a brand-new row for the SVM Regression search entry.
# ---------------------------------------------------------------------------
$ws.Rows("3:3").Insert()

# ---------------------------------------------------------------------------
# 4. Insert a new "Learning Type" column after the SMOTE Type column.
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Insert()

# ---------------------------------------------------------------------------
# 5. Row 2 -- WEKA Classification search, shortened model-file description
#    (set first so it claims the shared-string slot right after "na").
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "PerfSearch_Classification_TenfoldCV_WEKA_SMOTE217"

# ---------------------------------------------------------------------------
# 6. Header row -- new "Learning Type" column, and the "SVM Class." value
#    for the row above.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Learning Type"
$ws.Range("C2").Value = "SVM Class."

$ws.Range("G2").Value = 0.94930879999999995
$ws.Range("H2").Value = 0.94009220000000004
$ws.Range("I2").Value = 0.95852530000000002
$ws.Range("J2").Value = 0.89877019999999996
$ws.Range("K2").Value = 0.94700460829493105
$ws.Range("L2").Value = 0.93548387096774199
$ws.Range("M2").Value = 0.95852534562212
$ws.Range("N2").Value = 0.89424663014522598
$ws.Range("O2").Value = 0.9375
$ws.Range("P2").Value = 0.98039220000000005
$ws.Range("Q2").Value = 0.76923079999999999
$ws.Range("R2").Value = 0.79940219999999995

# ---------------------------------------------------------------------------
# 7. Row 3 -- new WEKA Regression search entry (metadata only, no results
#    yet -- this is the work moved to the "ForLater" folder).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "PerfSearch_Regression_TenfoldCV_WEKA_SMOTE217"
$ws.Range("B3").Value = "217+217"
$ws.Range("C3").Value = "SVM Regress."
$ws.Range("D3").Value = "1500-2800/50"
$ws.Range("A3").WrapText = $true

# ---------------------------------------------------------------------------
# 8. Row 5 ("Ahmed et al.") and row 6 ("Yang et al.") only moved one column
#    to the right because of the new column; their numbers are left exactly
#    as they already are (re-assigning them would lose the original clean
#    decimal representation).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 9. Row 6 -- "na" marker moves from column E to column F as text (it used
#    to sit one column closer to A before the new column was inserted).
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = "na"

# ---------------------------------------------------------------------------
# 10. Restore frozen header pane / selection state.
# ---------------------------------------------------------------------------
$ws.Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E3").Select()
